# "Add files via upload" — re-upload of the FortiNet model-comparison
# workbook with a couple of small manual fixes made in Excel before the
# file was saved:
#
#   1. A handful of rows got a stray "." typed into column A (outside the
#      Table1 data range, which still starts at column B). This shows up
#      as new A10:A15 / A24:A25 cells and pushes the sheet's used range
#      from B2:O28 out to A2:O28.
#   2. The "Firewall Latency" row header had a typo in its rich-text
#      label — the Symbol-font character meant to render as the Greek
#      "micro" glyph (µs) was being displayed literally, so it was
#      reworded to spell out "Firewall Latency (microseconds)" instead,
#      keeping the same Symbol-font character in the middle of the word.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Stray "." markers added in column A for these rows.
$dotRows = 10,11,12,13,14,15,24,25
foreach ($r in $dotRows) {
    $ws.Range("A$r").Value = "."
}

# 2) Reword the "Firewall Latency" label (row 10, column B).
#    Runs:  "Firewall Latency (micro" | " " (Symbol font) | "seconds)"
$ws.Range("B10").Value = "Firewall Latency (micro seconds)"
$ws.Range("B10").Characters(24, 1).Font.Name = "Symbol"
$ws.Range("B10").Characters(25, 9).Font.Name = "Aptos Narrow"

# Leave the selection on the cell that was last edited, as in the saved file.
$ws.Range("B10").Select()
